# Thêm cột "Trạng thái đại lý" vào báo cáo tổng hợp bán hàng.
# Cột mới được chèn trước cột D hiện tại (Đại lý mua), đẩy toàn bộ các cột
# từ D trở đi sang phải một vị trí.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chèn một cột mới tại vị trí D - Excel tự động dịch chuyển mergeCells,
# dimension, row spans và kế thừa định dạng từ cột bên trái (C).
$ws.Columns("D").Insert()

# Độ rộng cột mới (~20.29 ký tự trong bản gốc; gần nhất có thể đặt qua
# ColumnWidth trong môi trường này).
$ws.Columns("D").ColumnWidth = 19.5

# Tiêu đề cột (hàng 7) và ô dữ liệu mẫu (hàng 9) cho cột mới.
$ws.Range("D7").Value = "Trạng thái đại lý"
$ws.Range("D9").Value = "{{ReportSalesOrderGenerals.SalesOrders.BuyerStoreStatusName}}"

# Chiều cao hàng 9 giảm từ 195 xuống 165 do có thêm cột.
$ws.Rows(9).RowHeight = 165

# Vùng chọn hiện tại của sheet chuyển sang D13.
$ws.Range("D13").Select() | Out-Null
